# Bug fix for calculating achieved capacity factors for resources with no
# prior capacity. We now use calculated hypothetical capacity factors in
# every year except the start year for least cost dispatch resources to
# avoid data discontinuities.
#
# Concretely: on the "CSC-CSCSoCECBiaSY" sheet (Capacity Supply Curve ->
# Share of Cost Effective Capacity Built in a Single Year / "Max share of
# existing capacity buildable"), raise the max share of existing capacity
# buildable from 0.5 to 0.55 for every resource row (rows 2-25) across all
# year columns (B:AE, years 2021-2050).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")

$ws.Range("B2:AE25").Value = 0.55
